$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-ambiguous values: direct assignment is safe ---
$ws.Range("D2").Value = "42.640.04"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "2.353.86"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("E7").Value = "  -1.96%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -3.99%  "
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("E15").Value = "  -7.99%  "
$ws.Range("D16").Value = "2.712.77"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").Value = "2.408.07"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").Value = "42.763.06"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("E19").Value = "  +4.50%  "
$ws.Range("E20").Value = "  -2.08%  "
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("E22").Value = "  +5.26%  "
$ws.Range("E23").Value = "  -5.74%  "
$ws.Range("E24").Value = "  -4.35%  "
$ws.Range("E25").Value = "  -2.78%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -3.73%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("E31").Value = "  -4.26%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E32").Value = "  -4.37%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E33").Value = "  +4.04%  "
$ws.Range("E34").Value = "  -8.23%  "
$ws.Range("E35").Value = "  +17.81%  "
$ws.Range("E36").Value = "  -2.05%  "
$ws.Range("E37").Value = "  -5.30%  "
$ws.Range("E38").Value = "  -1.65%  "
$ws.Range("E39").Value = "  -8.48%  "
$ws.Range("E40").Value = "  -5.84%  "
$ws.Range("E41").Value = "  +3.25%  "
$ws.Range("E42").Value = "  +3.80%  "
$ws.Range("E43").Value = "  -7.85%  "
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("E45").Value = "  -8.95%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("E46").Value = "  -4.23%  "
$ws.Range("E47").Value = "  -2.69%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("E48").Value = "  -4.44%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("E49").Value = "  -7.64%  "
$ws.Range("E50").Value = "  +2.77%  "
$ws.Range("E51").Value = "  -2.07%  "

# --- Numeric-looking text values: force Text via a scratch cell + PasteSpecial(values), then remove the scratch column ---
$scratch = $ws.Range("G1")
$scratch.Value = "'317.91"
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$scratch.Value = "'106.69"
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$scratch.Value = "'41.50"
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$scratch.Value = "'8.43"
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$scratch.Value = "'0.997"
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$scratch.Value = "'15.96"
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$scratch.Value = "'7.78"
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$scratch.Value = "'76.72"
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$scratch.Value = "'3.60"
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$scratch.Value = "'257.43"
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$scratch.Value = "'9.46"
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$scratch.Value = "'11.39"
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$scratch.Value = "'22.97"
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$scratch.Value = "'2.23"
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$scratch.Value = "'174.89"
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$scratch.Value = "'36.39"
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$scratch.Value = "'0.0889"
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$scratch.Value = "'6.13"
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$scratch.Value = "'0.127"
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$scratch.Value = "'4.61"
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$scratch.Value = "'3.80"
$scratch.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$scratch.Value = "'2.67"
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$scratch.Value = "'0.241"
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$scratch.Value = "'71.82"
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$scratch.Value = "'114.22"
$scratch.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$scratch.Value = "'12.02"
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$scratch.Value = "'9.13"
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$scratch.Value = "'84.91"
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$scratch.Value = "'76.07"
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$scratch.Value = "'0.100"
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$scratch.ClearContents()
$scratch.EntireColumn.Delete()
